# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect the latest scraped totals (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 631
$ws1.Range("F4").Value = 659
$ws1.Range("F5").Value = 570
$ws1.Range("F6").Value = 315
$ws1.Range("F7").Value = 2805
$ws1.Range("F9").Value = 7907
$ws1.Range("F11").Value = 474
$ws1.Range("F13").Value = 367
$ws1.Range("F14").Value = 50

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 631
$ws4.Range("F4").Value = 659
$ws4.Range("F5").Value = 570
$ws4.Range("F6").Value = 315
$ws4.Range("F9").Value = 2805
$ws4.Range("F11").Value = 7907
$ws4.Range("F13").Value = 474
$ws4.Range("F17").Value = 367
$ws4.Range("F18").Value = 50
